# Economic Dashboard update - 2026-01-14
# Applies the cell value / highlight-style changes described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Style changes: move the "recently updated" yellow highlight (style used on
# N22, fillId none == style 48) off of N5/N10/N11/N12, and onto C17/C18/C40/
# C41 (style used on C24, yellow fill == style 49). We reuse existing
# identically-formatted donor cells via Copy + PasteSpecial(xlPasteFormats)
# so the workbook's existing style indices are reused rather than new
# duplicate styles being minted.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

# Remove highlight (s=49 -> s=48) on N5, N10, N11, N12
$ws.Range("N22").Copy()
$ws.Range("N5").PasteSpecial($xlPasteFormats)
$ws.Range("N10").PasteSpecial($xlPasteFormats)
$ws.Range("N11").PasteSpecial($xlPasteFormats)
$ws.Range("N12").PasteSpecial($xlPasteFormats)

# Add highlight (s=48 -> s=49) on C17, C18, C40, C41
$ws.Range("C24").Copy()
$ws.Range("C17").PasteSpecial($xlPasteFormats)
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("C40").PasteSpecial($xlPasteFormats)
$ws.Range("C41").PasteSpecial($xlPasteFormats)

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 5 - ADP, Total NonFarm Private: date unchanged, just dropped highlight
# (handled above). No value changes for row 5.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Rows 10-12 - JOLTS Openings/Hires/Separations Rate: date unchanged, just
# dropped highlight (handled above). No value changes for rows 10-12.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Row 17 - Retail Sales, M/M % Delta: new latest date + refreshed lag series
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 45962
$ws.Range("F17").Value = 0.006141545394387826
$ws.Range("G17").Value = -0.001065294349023249
$ws.Range("H17").Value = 0.0006724067240673204
$ws.Range("I17").Value = 0.00545946488174831
$ws.Range("J17").Value = 0.006492096487988874

# ---------------------------------------------------------------------------
# Row 18 - Retail Sales, Y/Y % Delta: new latest date + refreshed lag series
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 45962
$ws.Range("F18").Value = 0.03331905781584583
$ws.Range("G18").Value = 0.03264482062432054
$ws.Range("H18").Value = 0.04144341481107452
$ws.Range("I18").Value = 0.04972605550048132
$ws.Range("J18").Value = 0.04134309243240536

# ---------------------------------------------------------------------------
# Row 29 - 5yr, 5yr Forward: new latest date + refreshed lag series
# ---------------------------------------------------------------------------
$ws.Range("N29").Value = 46035
$ws.Range("Q29").Value = 2.23
$ws.Range("R29").Value = 2.22
$ws.Range("S29").Value = 2.24
$ws.Range("T29").Value = 2.23

# ---------------------------------------------------------------------------
# Row 30 - 10yr TIPS: new latest date + refreshed lag series
# ---------------------------------------------------------------------------
$ws.Range("N30").Value = 46035
$ws.Range("Q30").Value = 2.3
$ws.Range("R30").Value = 2.29
$ws.Range("S30").Value = 2.28

# ---------------------------------------------------------------------------
# Row 40 - New Home Sales, SAAR (Thousands): new latest date + refreshed
# lag series
# ---------------------------------------------------------------------------
$ws.Range("C40").Value = 45931
$ws.Range("F40").Value = 737
$ws.Range("G40").Value = 738
$ws.Range("H40").Value = 711
$ws.Range("I40").Value = 639
$ws.Range("J40").Value = 662

# ---------------------------------------------------------------------------
# Row 41 - New Home Sales, Y/Y % Delta: new latest date + refreshed lag
# series
# ---------------------------------------------------------------------------
$ws.Range("C41").Value = 45931
$ws.Range("F41").Value = 0.1867954911433172
$ws.Range("G41").Value = 0.02928870292887029
$ws.Range("H41").Value = 0.02597402597402598
$ws.Range("I41").Value = -0.1
$ws.Range("J41").Value = -0.01341281669150522

# ---------------------------------------------------------------------------
# Row 47 - FFR: new latest date only
# ---------------------------------------------------------------------------
$ws.Range("N47").Value = 46034

# ---------------------------------------------------------------------------
# Row 48 - 2y UST: new latest date + refreshed lag series
# ---------------------------------------------------------------------------
$ws.Range("N48").Value = 46034
$ws.Range("R48").Value = 3.54
$ws.Range("S48").Value = 3.49
$ws.Range("U48").Value = 3.47

# ---------------------------------------------------------------------------
# Row 49 - 5y UST: new latest date + refreshed lag series
# ---------------------------------------------------------------------------
$ws.Range("N49").Value = 46034
$ws.Range("Q49").Value = 3.77
$ws.Range("R49").Value = 3.75
$ws.Range("S49").Value = 3.74
$ws.Range("T49").Value = 3.7
$ws.Range("U49").Value = 3.72

# ---------------------------------------------------------------------------
# Row 50 - 10y UST: new latest date + refreshed lag series
# ---------------------------------------------------------------------------
$ws.Range("N50").Value = 46034
$ws.Range("Q50").Value = 4.19
$ws.Range("R50").Value = 4.18
$ws.Range("S50").Value = 4.19
$ws.Range("T50").Value = 4.15
$ws.Range("U50").Value = 4.18

# ---------------------------------------------------------------------------
# Row 52 - BAA: new latest date + refreshed lag series
# ---------------------------------------------------------------------------
$ws.Range("N52").Value = 46034
$ws.Range("Q52").Value = 5.89
$ws.Range("R52").Value = 5.88
$ws.Range("S52").Value = 5.92
$ws.Range("T52").Value = 5.88
